# Refactoring; add tests for all preflight.py methods
#
# The error-code table on the "Errors_" sheet (internal codeName "Sheet2")
# is renumbered: every iCode value in rows 9-36 (column A) is shifted up by
# 70 (e.g. 130 -> 200, 131 -> 201, ... 404 -> 474). The rest of each row
# (Class / Locn / Msg_String in columns B-D) is unchanged.
#
# The sheet selection is also moved from the old editing location (D64) to
# the block of cells that was just edited (E9:E36), with the active cell at
# the top of that block (E9).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errors_")

# Old iCode -> new iCode for rows 9 through 36 of column A.
$newCodes = @(200, 201, 220, 222, 240, 243, 260, 264, 280, 285, 300, 306, `
              320, 327, 340, 348, 360, 369, 380, 390, 400, 411, 420, 432, `
              440, 453, 460, 474)

$startRow = 9
for ($i = 0; $i -lt $newCodes.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newCodes[$i]
}

# Move the selection to reflect the range that was just edited.
$ws.Activate()
$ws.Range("E9:E36").Select()
